# Memphis Grizzlies stats workbook update:
#   - Insert two new game-log sheets "Rebounds" and "3PM" right after "Assists"
#     (pushing "Avg Points" / "Avg Assists" later in the tab order).
#   - Append two new summary sheets "Avg Rebounds" and "Avg 3PM" at the end.

$wb = $excel.ActiveWorkbook

$gameHeaders = @(
    "Game Time (PST)","Opponent","Jaylen Wells","Kentavious Caldwell-Pope",
    "Vince Williams Jr.","Santi Aldama","Jaren Jackson Jr.","Javon Small",
    "Ja Morant","PJ Hall","Olivier-Maxence Prosper","Cedric Coward",
    "Cam Spencer","Jock Landale","GG Jackson","John Konchar"
)

$REBOUNDS_DATA = @(
    @("2025-10-22","NOP",3,4,0,6,8,0,3,0,5,3,0,4,0,3),
    @("2025-10-24","MIA",4,3,0,8,3,2,2,3,4,2,2,5,2,1),
    @("2025-10-25","IND",3,4,0,8,2,2,3,0,2,6,6,8,0,5),
    @("2025-10-27","GSW",1,2,0,6,7,3,3,1,0,4,2,6,0,4),
    @("2025-10-29","PHX",3,2,0,10,6,0,8,0,0,4,2,3,0,3),
    @("2025-10-31","LAL",7,2,2,10,4,0,1,0,0,10,0,5,0,1),
    @("2025-11-02","TOR",5,0,3,4,9,1,0,0,2,5,4,6,0,2),
    @("2025-11-03","DET",6,4,2,5,4,0,5,0,0,6,2,5,0,0),
    @("2025-11-05","HOU",3,3,4,4,6,0,5,2,1,9,2,3,0,0),
    @("2025-11-07","DAL",5,3,3,2,4,0,5,0,0,9,2,7,0,2),
    @("2025-11-09","OKC",1,4,6,9,7,0,3,0,0,10,1,5,1,3)
)

$THREEPM_DATA = @(
    @("2025-10-22","NOP",4,2,0,0,0,0,2,0,0,0,2,1,0,0),
    @("2025-10-24","MIA",3,2,0,1,1,2,0,1,0,2,0,0,0,1),
    @("2025-10-25","IND",0,1,0,1,1,3,0,0,0,6,2,1,0,0),
    @("2025-10-27","GSW",0,1,0,2,2,0,0,0,0,1,3,2,0,1),
    @("2025-10-29","PHX",1,1,0,2,2,1,3,0,0,2,2,0,0,1),
    @("2025-10-31","LAL",2,2,1,1,2,0,0,0,0,0,2,2,0,2),
    @("2025-11-02","TOR",1,2,2,2,1,0,0,0,0,0,1,2,0,0),
    @("2025-11-03","DET",0,0,0,1,4,0,0,0,0,3,4,2,0,0),
    @("2025-11-05","HOU",1,0,2,1,0,0,1,0,0,0,5,3,0,0),
    @("2025-11-07","DAL",1,2,0,3,2,0,0,0,0,3,2,0,0,0),
    @("2025-11-09","OKC",2,1,1,1,3,0,2,0,0,3,1,2,0,1)
)

$AVG_REBOUNDS_DATA = @(
    @("Santi Aldama",6.545454545454546),
    @("Cedric Coward",6.181818181818182),
    @("Jaren Jackson Jr.",5.454545454545454),
    @("Jock Landale",5.181818181818182),
    @("Ja Morant",3.8),
    @("Jaylen Wells",3.727272727272727),
    @("Vince Williams Jr.",3.333333333333333),
    @("Kentavious Caldwell-Pope",2.818181818181818),
    @("John Konchar",2.666666666666667),
    @("Cam Spencer",2.090909090909091),
    @("Olivier-Maxence Prosper",1.75),
    @("Javon Small",1),
    @("PJ Hall",1),
    @("GG Jackson",0.6)
)

$AVG_3PM_DATA = @(
    @("Cam Spencer",2.181818181818182),
    @("Cedric Coward",1.818181818181818),
    @("Jaren Jackson Jr.",1.636363636363636),
    @("Jaylen Wells",1.363636363636364),
    @("Santi Aldama",1.363636363636364),
    @("Jock Landale",1.363636363636364),
    @("Kentavious Caldwell-Pope",1.272727272727273),
    @("Vince Williams Jr.",1),
    @("Ja Morant",0.8),
    @("Javon Small",0.75),
    @("John Konchar",0.6666666666666666),
    @("PJ Hall",0.1666666666666667),
    @("Olivier-Maxence Prosper",0),
    @("GG Jackson",0)
)

function Write-GameLogSheet($sheet, $headers, $data) {
    for ($c = 0; $c -lt $headers.Length; $c++) {
        $cell = $sheet.Cells.Item(1, $c + 1)
        $cell.Value = $headers[$c]
        $cell.Font.Bold = $true
        $cell.HorizontalAlignment = -4108
        $cell.VerticalAlignment = -4160
        $cell.Borders.LineStyle = 1
    }
    for ($r = 0; $r -lt $data.Length; $r++) {
        $row = $data[$r]
        for ($c = 0; $c -lt $row.Length; $c++) {
            $cell = $sheet.Cells.Item($r + 2, $c + 1)
            if ($c -eq 0) {
                $cell.NumberFormat = "@"
            }
            $cell.Value = $row[$c]
        }
    }
}

function Write-AvgSheet($sheet, $data) {
    $h1 = $sheet.Cells.Item(1, 1)
    $h1.Value = "Player"
    $h1.Font.Bold = $true
    $h1.HorizontalAlignment = -4108
    $h1.VerticalAlignment = -4160
    $h1.Borders.LineStyle = 1

    $h2 = $sheet.Cells.Item(1, 2)
    $h2.Value = $sheet.Name
    $h2.Font.Bold = $true
    $h2.HorizontalAlignment = -4108
    $h2.VerticalAlignment = -4160
    $h2.Borders.LineStyle = 1

    for ($r = 0; $r -lt $data.Length; $r++) {
        $row = $data[$r]
        $sheet.Cells.Item($r + 2, 1).Value = $row[0]
        $sheet.Cells.Item($r + 2, 2).Value = $row[1]
    }
}

# --- Insert "Rebounds" and "3PM" game-log sheets right after "Assists" ---
$assists = $wb.Worksheets.Item("Assists")

$reboundsSheet = $wb.Worksheets.Add($null, $assists)
$reboundsSheet.Name = "Rebounds"
Write-GameLogSheet $reboundsSheet $gameHeaders $REBOUNDS_DATA

$threePmSheet = $wb.Worksheets.Add($null, $reboundsSheet)
$threePmSheet.Name = "3PM"
Write-GameLogSheet $threePmSheet $gameHeaders $THREEPM_DATA

# --- Append "Avg Rebounds" and "Avg 3PM" summary sheets at the end ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

$avgReboundsSheet = $wb.Worksheets.Add($null, $lastSheet)
$avgReboundsSheet.Name = "Avg Rebounds"
Write-AvgSheet $avgReboundsSheet $AVG_REBOUNDS_DATA

$avg3pmSheet = $wb.Worksheets.Add($null, $avgReboundsSheet)
$avg3pmSheet.Name = "Avg 3PM"
Write-AvgSheet $avg3pmSheet $AVG_3PM_DATA
